# Actualización automática 2025-07-10 17:25:09
#
# Registers the "julio" (July) sales for two clients of
# ALMEIDA CUATIN JHONATHANN CARLOS (AUCANSHALA ALLAICA FREDDY HERNAN and
# PALATE CHUCARALAO JOSE ISRAEL) in the "240X80 PORCELANATO" product
# group, and ripples that change through every dependent total /
# percentage cell on the three report sheets.
#
# NOTE: this runtime's Range/Cells ".Value" getter does not reliably
# surface numeric values when read back inside an expression, so every
# read below uses ".Value2" instead (writes work fine with either).

$wb = $excel.ActiveWorkbook

$wsGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl   = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": column D = "240X80 PORCELANATO"
# ---------------------------------------------------------------
$wsGrupo.Cells.Item(4, 4).Value2 = 950.4     # AUCANSHALA ALLAICA FREDDY HERNAN
$wsGrupo.Cells.Item(21, 4).Value2 = 570.24   # PALATE CHUCARALAO JOSE ISRAEL

# Totals row (row 32) shows, per product column, how many of the 30
# clients (rows 2-31) have a non-zero figure, as "n de 30". Only column
# D (index 4) is affected by the values changed above.
$countD = 0
for ($row = 2; $row -le 31; $row++) {
    if ($wsGrupo.Cells.Item($row, 4).Value2 -ne 0) { $countD++ }
}
$wsGrupo.Cells.Item(32, 4).Value2 = "$countD de 30"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": column F = "julio"
# ---------------------------------------------------------------
$oldF4  = $wsMensual.Cells.Item(4, 6).Value2
$oldF21 = $wsMensual.Cells.Item(21, 6).Value2
$newF4  = 950.4
$newF21 = 570.24

$wsMensual.Cells.Item(4, 6).Value2 = $newF4
$wsMensual.Cells.Item(21, 6).Value2 = $newF21

# Totals row (row 32) = SUM of column F across rows 2-31; adjust the
# existing total by the same delta rather than re-summing every row,
# so unrelated totals (C32/D32/E32/G32) keep their original values.
$deltaF = ($newF4 - $oldF4) + ($newF21 - $oldF21)
$wsMensual.Cells.Item(32, 6).Value2 = $wsMensual.Cells.Item(32, 6).Value2 + $deltaF

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": rows 2 (240X80 PORCELANATO) and 18 (TOTAL)
# C = PRESUPUESTO, D = VENTA, E = POR CUMPLIR (C-D), F = CUMPLIMIENTO (D/C)
# ---------------------------------------------------------------

# Column D ("VENTA") is a touch wider now that it holds bigger numbers.
# (ColumnWidth is stored with the usual +5/6 char padding offset, so
# feed it 13 - 5/6 to land exactly on a stored width of 13.)
$wsCumpl.Columns.Item(4).ColumnWidth = 12.166666666666666

$oldD2 = $wsCumpl.Cells.Item(2, 4).Value2
$newD2 = 1520.64
$wsCumpl.Cells.Item(2, 4).Value2 = $newD2

$c2 = $wsCumpl.Cells.Item(2, 3).Value2
$wsCumpl.Cells.Item(2, 5).Value2 = $c2 - $newD2
$wsCumpl.Cells.Item(2, 6).Value2 = $newD2 / $c2

# Row 18 (TOTAL): only D/E/F move, driven by the D2 delta above.
$deltaD = $newD2 - $oldD2
$oldD18 = $wsCumpl.Cells.Item(18, 4).Value2
$newD18 = $oldD18 + $deltaD
$wsCumpl.Cells.Item(18, 4).Value2 = $newD18

$c18 = $wsCumpl.Cells.Item(18, 3).Value2
$wsCumpl.Cells.Item(18, 5).Value2 = $c18 - $newD18
$wsCumpl.Cells.Item(18, 6).Value2 = $newD18 / $c18
